$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: add the 0.5h estimate for the new Q (week) column ---
$ws.Range("Q22").Value = 0.5

# --- Row 23: fill in hours spent across N/O/Q columns ---
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 3
$ws.Range("Q23").Value = 1

# --- Row 29 (item 22): new task "Realizar BBDD + inserts" ---
$ws.Range("C29").Value = "Realizar BBDD + inserts"
$ws.Range("Q29").Value = 11
$ws.Range("AA29").Value = "Realizado el 14/11/18"
$ws.Range("AB29").Value = 15

# --- Row 30 (item 23): new task "Listview resultados" ---
# Bring over the same look-and-feel (borders/fill/font) used by row 29's
# label cells, then fill in the content.
$ws.Range("C29:I29").Copy()
$ws.Range("C30:I30").PasteSpecial(-4122)
$ws.Range("C30").Value = "Listview resultados"
$ws.Range("Q30").Value = 4

# --- Move the active selection like the author left it ---
$ws.Range("Q31").Select()
